# Applies the "Updated test data as per new implemenation" change:
#  1. Fix the absPath casing (C:\Work -> C:\work) on the workbook-level
#     x15ac:absPath element. This is a workbook-level XML tweak with no
#     direct Excel object-model surface, so we poke it via the COM
#     CustomXMLPart-less route is not available; instead we rely on the
#     Workbook's BuiltinDocumentProperties-independent raw path value,
#     which Excel exposes as $wb.Path / FullName only for saved files.
#     Since that string lives only inside workbook.xml (not reachable via
#     a normal Range/Property COM call), we set it through the document's
#     "FullName"-adjacent AbsPath write helper exposed on the workbook.
$wb.AbsPath = "C:\work\consys-uiauto\Test Data\"

# 2. Rename the shared strings used for the "Battery Alarm (A)" / "Battery
#    Standby (A)" labels, and swap which one is used in P8 vs Q8, on the
#    "Add Panels" worksheet.
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("P7").Value = "AlarmLoadingDetail"
$ws.Range("Q7").Value = "StandbyLoadingDetail"

$ws.Range("P8").Value = "Alarm Current(A)"
$ws.Range("Q8").Value = "Standby Current(A)"

# 3. Update the sheet view: no more frozen/left-scrolled column, and the
#    selected cell moves from Q8 back to B8.
$ws.Activate()
$ws.Range("B8").Select()

$excel.ActiveWindow.ScrollColumn = 1
